$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 89, pushing the existing rows 89..186 down to 90..187.
$ws.Rows.Item(89).Insert()

# Populate the new row 89 with a new weekly data point (matches the other
# cells of the old row 89, which is now row 90, except for the columns below).
$ws.Cells.Item(89, 1).Value = 9
$ws.Cells.Item(89, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(89, 3).Value = "Metropolitana"
$ws.Cells.Item(89, 4).Value = 45175
$ws.Cells.Item(89, 5).Value = 13
$ws.Cells.Item(89, 6).Value = 100112022
$ws.Cells.Item(89, 7).Value = "Arveja Verde"
$ws.Cells.Item(89, 8).Value = "Perfection"
$ws.Cells.Item(89, 9).Value = "Primera"
$ws.Cells.Item(89, 10).Value = 52
$ws.Cells.Item(89, 11).Value = 26000
$ws.Cells.Item(89, 12).Value = 28000
$ws.Cells.Item(89, 13).Value = 27000
$ws.Cells.Item(89, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(89, 15).Value = "Provincia de Limar" + [char]0xED
$ws.Cells.Item(89, 16).Value = 1080
$ws.Cells.Item(89, 17).Value = 25
$ws.Cells.Item(89, 18).Value = "Hortaliza"
